$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G6").Value = "Prakruti Sinha"
